$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 127, shifting existing rows 127-196 down to 128-197.
$ws.Rows.Item(127).Insert()

# Populate the newly inserted row 127 with its data.
$ws.Cells.Item(127, 1).Value = 10
$ws.Cells.Item(127, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(127, 3).Value = "La Araucanía"
$ws.Cells.Item(127, 4).Value = 44603
$ws.Cells.Item(127, 5).Value = 9
$ws.Cells.Item(127, 6).Value = "Fruta"
$ws.Cells.Item(127, 7).Value = 100102
$ws.Cells.Item(127, 8).Value = "Cítricos"
$ws.Cells.Item(127, 9).Value = 100102006
$ws.Cells.Item(127, 10).Value = "Pomelo"
$ws.Cells.Item(127, 11).Value = "Red Blush"
$ws.Cells.Item(127, 12).Value = "Primera"
$ws.Cells.Item(127, 13).Value = 80
$ws.Cells.Item(127, 14).Value = 13000
$ws.Cells.Item(127, 15).Value = 13000
$ws.Cells.Item(127, 16).Value = 13000
$ws.Cells.Item(127, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(127, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(127, 19).Value = 867
$ws.Cells.Item(127, 20).Value = 15
